$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new predidx (D) and pred_name (E) values
$updates = @(
    @{ Row = 11;  D = "[1, 0, 0, 0, 1, 0, 0]"; E = "['Normal', 'RegulationViolation']" },
    @{ Row = 12;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 15;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 24;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 25;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 26;  D = "[0, 0, 0, 0, 0, 0, 1]"; E = "['SoftwareFault']" },
    @{ Row = 27;  D = "[0, 0, 1, 0, 0, 0, 1]"; E = "['HardwareFault', 'SoftwareFault']" },
    @{ Row = 39;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 56;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 58;  D = "[1, 0, 0, 1, 0, 0, 0]"; E = "['Normal', 'ParamViolation']" },
    @{ Row = 73;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 81;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 82;  D = "[1, 1, 1, 0, 0, 0, 0]"; E = "['Normal', 'SurroundingEnvironment', 'HardwareFault']" },
    @{ Row = 83;  D = "[1, 1, 0, 0, 0, 0, 0]"; E = "['Normal', 'SurroundingEnvironment']" },
    @{ Row = 89;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 113; D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" },
    @{ Row = 116; D = "[1, 0, 0, 0, 0, 1, 0]"; E = "['Normal', 'CommunicationIssue']" }
)

foreach ($u in $updates) {
    $ws.Range("D$($u.Row)").Value = $u.D
    $ws.Range("E$($u.Row)").Value = $u.E
}
